# The data in this sheet got re-shuffled across rows 2,3,5,6,7,8,9,10,11,12,13
# (row 4 stays put). Only the columns D (Fecha), L (Calidad), M (Volumen),
# N (Precio minimo), O (Precio maximo), P (Precio promedio ponderado),
# Q (Unidad de comercializacion), R (Origen), S (Precio $/Kg) and
# T (Kg / unidad) are affected for each row - the mapping below says which
# *old* row's values end up in which *new* row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns (by letter) that carry the values being shuffled between rows.
$cols = @("D", "L", "M", "N", "O", "P", "Q", "R", "S", "T")

# newRow -> oldRow (i.e. newRow should end up holding what oldRow used to hold)
$mapping = @{
    2  = 7
    3  = 8
    5  = 12
    6  = 11
    7  = 2
    8  = 3
    9  = 13
    10 = 5
    11 = 6
    12 = 9
    13 = 10
}

# Snapshot the "before" values for every row/column involved, so that
# overwriting one row doesn't clobber the source data needed by another.
$snapshot = @{}
foreach ($row in ($mapping.Values | Select-Object -Unique)) {
    $rowData = @{}
    foreach ($col in $cols) {
        $rowData[$col] = $ws.Range("$col$row").Value2
    }
    $snapshot[$row] = $rowData
}

# Now write the shuffled values back into the new rows.
foreach ($newRow in $mapping.Keys) {
    $oldRow = $mapping[$newRow]
    $rowData = $snapshot[$oldRow]
    foreach ($col in $cols) {
        $ws.Range("$col$newRow").Value2 = $rowData[$col]
    }
}
